$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 fixes ---
# B14 was stored as text "3"; it should become a genuine number 3.
$ws.Range("B14").Value = 3
# C14 was stored as text "nan"; it should become blank.
$ws.Range("C14").Value = ""

# --- New row 15 ---
$ws.Range("A15").Value = "parisk"

# B15 must stay text "4" (not be auto-converted to the number 4), so force
# the cell to a Text format before assigning the numeric-looking string.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "4"

$ws.Range("C15").Value = "well carried through"
$ws.Range("D15").Value = "FBK"
$ws.Range("E15").Value = "EXP"
$ws.Range("F15").Value = "ea04c829-c996-4167-8585-03efb193cd41"
$ws.Range("G15").Value = "ByOExmWAb_annotated.xlsx"
$ws.Range("H15").Value = "The experiments were well carried through and very thorough."
$ws.Range("I15").Value = "Needs Revision"
